# adds printing, json rendering, smaller fixes
#
# Appends three new translation-key rows to the "translations" sheet:
#   56: source_wikidata / aus Wikidata
#   57: json / JSON / JSON / JSON   (same label across all three languages)
#   58: data / data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = "source_wikidata"
$ws.Range("B56").Value = "aus Wikidata"

$ws.Range("A57").Value = "json"
$ws.Range("B57").Value = "JSON"
$ws.Range("C57").Value = "JSON"
$ws.Range("D57").Value = "JSON"

$ws.Range("A58").Value = "data"
$ws.Range("B58").Value = "data"

# Move / record the active selection the way the source file shows it
# (bottom pane, one row below the newly-added data).
$ws.Range("A59").Select() | Out-Null
